$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(2, 8).Value = 581.61536
$ws_ALC.Cells.Item(2, 9).Value = 365.85715
$ws_ALC.Cells.Item(2, 11).Value = 365.85715
$ws_ALC.Cells.Item(2, 13).Value = -252.85715
$ws_ALC.Cells.Item(9, 8).Value = 267.57144
$ws_ALC.Cells.Item(9, 9).Value = 325.25
$ws_ALC.Cells.Item(9, 10).Value = 190.66667
$ws_ALC.Cells.Item(9, 11).Value = 325.25
$ws_ALC.Cells.Item(9, 12).Value = 190.66667
$ws_ALC.Cells.Item(9, 13).Value = -156.25
$ws_ALC.Cells.Item(9, 14).Value = -528.6666700000001
$ws_ALC.Cells.Item(18, 8).Value = 2718.8
$ws_ALC.Cells.Item(18, 9).Value = 2718.8
$ws_ALC.Cells.Item(18, 11).Value = 2718.8
$ws_ALC.Cells.Item(18, 13).Value = -2434.8
$ws_ALC.Cells.Item(33, 8).Value = 33334294
$ws_ALC.Cells.Item(33, 9).Value = 708.52
$ws_ALC.Cells.Item(33, 11).Value = 708.52
$ws_ALC.Cells.Item(33, 13).Value = -479.52
$ws_ALC.Cells.Item(38, 8).Value = 1750.16
$ws_ALC.Cells.Item(38, 9).Value = 264.1111
$ws_ALC.Cells.Item(38, 10).Value = 5571.4287
$ws_ALC.Cells.Item(38, 11).Value = 792.3333
$ws_ALC.Cells.Item(38, 12).Value = 16714.2861
$ws_ALC.Cells.Item(38, 13).Value = -420.3333
$ws_ALC.Cells.Item(38, 14).Value = -17458.2861
$ws_ALC.Cells.Item(47, 8).Value = 4800
$ws_ALC.Cells.Item(47, 9).Value = 0
$ws_ALC.Cells.Item(47, 10).Value = 4800
$ws_ALC.Cells.Item(47, 11).Value = 0
$ws_ALC.Cells.Item(47, 12).Value = 4800
$ws_ALC.Cells.Item(47, 13).ClearContents()
$ws_ALC.Cells.Item(47, 14).Value = -6744
$ws_ALC.Cells.Item(58, 8).Value = 21593.334
$ws_ALC.Cells.Item(58, 9).Value = 44131.43
$ws_ALC.Cells.Item(58, 10).Value = 1872.5
$ws_ALC.Cells.Item(58, 11).Value = 132394.29
$ws_ALC.Cells.Item(58, 12).Value = 5617.5
$ws_ALC.Cells.Item(58, 13).Value = -132244.29
$ws_ALC.Cells.Item(58, 14).Value = -5917.5
$ws_ALC.Cells.Item(113, 8).Value = 3480.0303
$ws_ALC.Cells.Item(113, 9).Value = 2866.7646
$ws_ALC.Cells.Item(113, 10).Value = 4131.625
$ws_ALC.Cells.Item(113, 11).Value = 2866.7646
$ws_ALC.Cells.Item(113, 12).Value = 4131.625
$ws_ALC.Cells.Item(113, 13).Value = 387.2354
$ws_ALC.Cells.Item(113, 14).Value = -10639.625

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(32, 8).Value = 4632.2964
$ws_ARM.Cells.Item(32, 9).Value = 3104.9592
$ws_ARM.Cells.Item(32, 11).Value = 3104.9592
$ws_ARM.Cells.Item(32, 13).Value = -2817.9592
$ws_ARM.Cells.Item(61, 8).Value = 1271.909
$ws_ARM.Cells.Item(61, 9).Value = 1078.6154
$ws_ARM.Cells.Item(61, 10).Value = 1551.1111
$ws_ARM.Cells.Item(61, 11).Value = 1078.6154
$ws_ARM.Cells.Item(61, 12).Value = 1551.1111
$ws_ARM.Cells.Item(61, 13).Value = -866.6153999999999
$ws_ARM.Cells.Item(61, 14).Value = -1975.1111
$ws_ARM.Cells.Item(95, 8).Value = 26166.4
$ws_ARM.Cells.Item(95, 10).Value = 26166.4
$ws_ARM.Cells.Item(95, 12).Value = 26166.4
$ws_ARM.Cells.Item(95, 14).Value = -31658.4
$ws_ARM.Cells.Item(136, 8).Value = 1271.909
$ws_ARM.Cells.Item(136, 9).Value = 1078.6154
$ws_ARM.Cells.Item(136, 10).Value = 1551.1111
$ws_ARM.Cells.Item(136, 11).Value = 3235.8462
$ws_ARM.Cells.Item(136, 12).Value = 4653.3333
$ws_ARM.Cells.Item(136, 13).Value = -685.8462
$ws_ARM.Cells.Item(136, 14).Value = -9753.3333

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(22, 8).Value = 232.3077
$ws_BSM.Cells.Item(22, 9).Value = 215.45454
$ws_BSM.Cells.Item(22, 10).Value = 325
$ws_BSM.Cells.Item(22, 11).Value = 215.45454
$ws_BSM.Cells.Item(22, 12).Value = 325
$ws_BSM.Cells.Item(22, 13).Value = -42.45454000000001
$ws_BSM.Cells.Item(22, 14).Value = -671
$ws_BSM.Cells.Item(105, 8).Value = 2494.375
$ws_BSM.Cells.Item(105, 9).Value = 2471.2693
$ws_BSM.Cells.Item(105, 10).Value = 2537.2856
$ws_BSM.Cells.Item(105, 11).Value = 2471.2693
$ws_BSM.Cells.Item(105, 12).Value = 2537.2856
$ws_BSM.Cells.Item(105, 13).Value = -724.2692999999999
$ws_BSM.Cells.Item(105, 14).Value = -6031.2856

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(4, 8).Value = 2984.3076
$ws_CRP.Cells.Item(4, 10).Value = 2984.3076
$ws_CRP.Cells.Item(4, 12).Value = 2984.3076
$ws_CRP.Cells.Item(4, 14).Value = -3208.3076

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(43, 8).Value = 1000
$ws_CUL.Cells.Item(43, 10).Value = 1000
$ws_CUL.Cells.Item(43, 12).Value = 3000
$ws_CUL.Cells.Item(43, 14).Value = -3228
$ws_CUL.Cells.Item(52, 8).Value = 655.5714
$ws_CUL.Cells.Item(52, 10).Value = 655.5714
$ws_CUL.Cells.Item(52, 12).Value = 1966.7142
$ws_CUL.Cells.Item(52, 14).Value = -2498.7142
$ws_CUL.Cells.Item(75, 8).Value = 4500
$ws_CUL.Cells.Item(75, 9).Value = 0
$ws_CUL.Cells.Item(75, 10).Value = 4500
$ws_CUL.Cells.Item(75, 11).Value = 0
$ws_CUL.Cells.Item(75, 12).Value = 13500
$ws_CUL.Cells.Item(75, 13).ClearContents()
$ws_CUL.Cells.Item(75, 14).Value = -15496
$ws_CUL.Cells.Item(78, 8).Value = 4500
$ws_CUL.Cells.Item(78, 9).Value = 0
$ws_CUL.Cells.Item(78, 10).Value = 4500
$ws_CUL.Cells.Item(78, 11).Value = 0
$ws_CUL.Cells.Item(78, 12).Value = 40500
$ws_CUL.Cells.Item(78, 13).ClearContents()
$ws_CUL.Cells.Item(78, 14).Value = -50484
$ws_CUL.Cells.Item(107, 8).Value = 563578.4
$ws_CUL.Cells.Item(107, 9).Value = 639.4783
$ws_CUL.Cells.Item(107, 10).Value = 803348.7
$ws_CUL.Cells.Item(107, 11).Value = 1918.4349
$ws_CUL.Cells.Item(107, 12).Value = 2410046.1
$ws_CUL.Cells.Item(107, 13).Value = 1.565100000000029
$ws_CUL.Cells.Item(107, 14).Value = -2413886.1
$ws_CUL.Cells.Item(131, 8).Value = 890.4697
$ws_CUL.Cells.Item(131, 10).Value = 994
$ws_CUL.Cells.Item(131, 12).Value = 2982
$ws_CUL.Cells.Item(131, 14).Value = -13062

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(116, 8).Value = 56742
$ws_GSM.Cells.Item(116, 9).Value = 0
$ws_GSM.Cells.Item(116, 10).Value = 56742
$ws_GSM.Cells.Item(116, 11).Value = 0
$ws_GSM.Cells.Item(116, 12).Value = 56742
$ws_GSM.Cells.Item(116, 13).ClearContents()
$ws_GSM.Cells.Item(116, 14).Value = -65920
$ws_GSM.Cells.Item(122, 8).Value = 1307.1666
$ws_GSM.Cells.Item(122, 9).Value = 1354
$ws_GSM.Cells.Item(122, 10).Value = 1166.6666
$ws_GSM.Cells.Item(122, 11).Value = 4062
$ws_GSM.Cells.Item(122, 12).Value = 3499.9998
$ws_GSM.Cells.Item(122, 13).Value = -1612
$ws_GSM.Cells.Item(122, 14).Value = -8399.9998

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(2, 8).Value = 125250
$ws_LTW.Cells.Item(2, 10).Value = 5888.8887
$ws_LTW.Cells.Item(2, 12).Value = 5888.8887
$ws_LTW.Cells.Item(2, 14).Value = -6112.8887
$ws_LTW.Cells.Item(13, 8).Value = 450.5
$ws_LTW.Cells.Item(13, 9).Value = 401
$ws_LTW.Cells.Item(13, 10).Value = 500
$ws_LTW.Cells.Item(13, 11).Value = 401
$ws_LTW.Cells.Item(13, 12).Value = 500
$ws_LTW.Cells.Item(13, 13).Value = -261
$ws_LTW.Cells.Item(13, 14).Value = -780
$ws_LTW.Cells.Item(55, 8).Value = 492.26315
$ws_LTW.Cells.Item(55, 9).Value = 220.6875
$ws_LTW.Cells.Item(55, 10).Value = 1940.6666
$ws_LTW.Cells.Item(55, 11).Value = 220.6875
$ws_LTW.Cells.Item(55, 12).Value = 1940.6666
$ws_LTW.Cells.Item(55, 13).Value = -47.6875
$ws_LTW.Cells.Item(55, 14).Value = -2286.6666
$ws_LTW.Cells.Item(132, 8).Value = 2525.8975
$ws_LTW.Cells.Item(132, 9).Value = 2422.5557
$ws_LTW.Cells.Item(132, 10).Value = 2666.818
$ws_LTW.Cells.Item(132, 11).Value = 7267.6671
$ws_LTW.Cells.Item(132, 12).Value = 8000.454000000001
$ws_LTW.Cells.Item(132, 13).Value = -4737.6671
$ws_LTW.Cells.Item(132, 14).Value = -13060.454

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(132, 8).Value = 6742.294
$ws_WVR.Cells.Item(132, 9).Value = 12030.143
$ws_WVR.Cells.Item(132, 10).Value = 3040.8
$ws_WVR.Cells.Item(132, 11).Value = 36090.429
$ws_WVR.Cells.Item(132, 12).Value = 9122.400000000001
$ws_WVR.Cells.Item(132, 13).Value = -33560.429
$ws_WVR.Cells.Item(132, 14).Value = -14182.4
